$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.16
$ws.Range("P2").Value = 4.6
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.42
$ws.Range("T2").Value = 3.25
$ws.Range("W2").Value = 18.5
$ws.Range("X2").Value = 35
$ws.Range("AB2").Value = 35
$ws.Range("AC2").Value = 9
$ws.Range("AD2").Value = 7.9
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 40
$ws.Range("AH2").Value = 10.25
$ws.Range("AM2").Value = 17
$ws.Range("AN2").Value = 7.1
$ws.Range("AP2").Value = 26
$ws.Range("AT2").Value = 3.25
$ws.Range("AU2").Value = 6.5
$ws.Range("AV2").Value = 45
$ws.Range("AX2").Value = 3.8
